# Trade #4 closed at 2026-02-17 13:33:31 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.93
$summary.Range("B4").Value = -0.07000000000000001
$summary.Range("B5").Value = -0.35
$summary.Range("B6").Value = 4
$summary.Range("B7").Value = 1
$summary.Range("B9").Value = 25

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.93000000000001
$status.Range("D4").Value = 4
$status.Range("E4").Value = -0.07000000000000001
$status.Range("F4").Value = -0.07000000000000001
$status.Range("G4").Value = 25

# --- New trade row data (Trade #4) ---
$tradeNum = 4
$tradeDate = "2026-02-17"
$tradeTime = "13:33:25"
$strategy = "MarketMaking"
$side = "DOWN"
$entryPrice = 0.91
$exitPrice = 0.93
$status4 = "CLOSED"
$pnlPct = 2.1978
$pnlDollar = 0.02
$capitalAfter = 99.93000000000001
$entrySlippage = 0
$exitSlippage = 0
$confidence = 0.6
$entryReason = "Normal spread capture: 19600 bps"
$exitReason = "early_exit"
$duration = 0.13

# --- All Trades sheet: append row 5 ---
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(5, 1).Value = $tradeNum
$allTrades.Cells.Item(5, 2).NumberFormat = "@"
$allTrades.Cells.Item(5, 2).Value = $tradeDate
$allTrades.Cells.Item(5, 2).Style = "Normal"
$allTrades.Cells.Item(5, 3).NumberFormat = "@"
$allTrades.Cells.Item(5, 3).Value = $tradeTime
$allTrades.Cells.Item(5, 3).Style = "Normal"
$allTrades.Cells.Item(5, 4).Value = $strategy
$allTrades.Cells.Item(5, 5).Value = $side
$allTrades.Cells.Item(5, 6).Value = $entryPrice
$allTrades.Cells.Item(5, 7).Value = $exitPrice
$allTrades.Cells.Item(5, 8).Value = $status4
$allTrades.Cells.Item(5, 9).Value = $pnlPct
$allTrades.Cells.Item(5, 10).Value = $pnlDollar
$allTrades.Cells.Item(5, 11).Value = $capitalAfter
$allTrades.Cells.Item(5, 12).Value = $entrySlippage
$allTrades.Cells.Item(5, 13).Value = $exitSlippage
$allTrades.Cells.Item(5, 14).Value = $confidence
$allTrades.Cells.Item(5, 15).Value = $entryReason
$allTrades.Cells.Item(5, 16).Value = $exitReason
$allTrades.Cells.Item(5, 17).Value = $duration

# --- MarketMaking sheet: append row 5 ---
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Cells.Item(5, 1).Value = $tradeNum
$marketMaking.Cells.Item(5, 2).NumberFormat = "@"
$marketMaking.Cells.Item(5, 2).Value = $tradeDate
$marketMaking.Cells.Item(5, 2).Style = "Normal"
$marketMaking.Cells.Item(5, 3).NumberFormat = "@"
$marketMaking.Cells.Item(5, 3).Value = $tradeTime
$marketMaking.Cells.Item(5, 3).Style = "Normal"
$marketMaking.Cells.Item(5, 4).Value = $strategy
$marketMaking.Cells.Item(5, 5).Value = $side
$marketMaking.Cells.Item(5, 6).Value = $entryPrice
$marketMaking.Cells.Item(5, 7).Value = $exitPrice
$marketMaking.Cells.Item(5, 8).Value = $status4
$marketMaking.Cells.Item(5, 9).Value = $pnlPct
$marketMaking.Cells.Item(5, 10).Value = $pnlDollar
$marketMaking.Cells.Item(5, 11).Value = $capitalAfter
$marketMaking.Cells.Item(5, 12).Value = $entrySlippage
$marketMaking.Cells.Item(5, 13).Value = $exitSlippage
$marketMaking.Cells.Item(5, 14).Value = $confidence
$marketMaking.Cells.Item(5, 15).Value = $entryReason
$marketMaking.Cells.Item(5, 16).Value = $exitReason
$marketMaking.Cells.Item(5, 17).Value = $duration
